$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: financial period headers (shift window forward one quarter) ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("E8").Value = "3 ماهه منتهی به 1400/03"
$ws.Range("F8").Value = "6 ماهه منتهی به 1400/06"
$ws.Range("G8").Value = "9 ماهه منتهی به 1400/09"
$ws.Range("H8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("I8").Value = "3 ماهه منتهی به 1401/03"
$ws.Range("J8").Value = "6 ماهه منتهی به 1401/06"
$ws.Range("K8").Value = "9 ماهه منتهی به 1401/09"
$ws.Range("L8").Value = "12 ماهه منتهی به 1401/12"
$ws.Range("M8").Value = "3 ماهه منتهی به 1402/03"

# --- Row 9: publish dates ---
$ws.Range("D9").Value = "1401-04-08 (8)"
$ws.Range("E9").Value = "1401-05-04 (3)"
$ws.Range("F9").Value = "1401-08-30 (4)"
$ws.Range("G9").Value = "1401-10-28 (2)"
$ws.Range("H9").Value = "1402-04-14 (9)"
$ws.Range("I9").Value = "1402-04-28 (3)"
$ws.Range("J9").Value = "1401-08-30 (2)"
$ws.Range("K9").Value = "1401-10-28"
$ws.Range("L9").Value = "1402-04-28 (3)"
$ws.Range("M9").Value = "1402-04-28"

# --- Data rows 11-27: shift cumulative figures forward one quarter + new quarter values ---
# Row 11
$ws.Range("D11").Value = 143234768
$ws.Range("E11").Value = 54684355
$ws.Range("F11").Value = 108851285
$ws.Range("G11").Value = 154394497
$ws.Range("H11").Value = 192628444
$ws.Range("I11").Value = 69507214
$ws.Range("J11").Value = 127578003
$ws.Range("K11").Value = 173935378
$ws.Range("L11").Value = 214213606
$ws.Range("M11").Value = 66850505
# Row 12
$ws.Range("D12").Value = -61344224
$ws.Range("E12").Value = -23861664
$ws.Range("F12").Value = -61769696
$ws.Range("G12").Value = -113390851
$ws.Range("H12").Value = -146246354
$ws.Range("I12").Value = -44141537
$ws.Range("J12").Value = -89353355
$ws.Range("K12").Value = -119387375
$ws.Range("L12").Value = -145108587
$ws.Range("M12").Value = -38739985
# Row 13
$ws.Range("D13").Value = 81890544
$ws.Range("E13").Value = 30822691
$ws.Range("F13").Value = 47081589
$ws.Range("G13").Value = 41003646
$ws.Range("H13").Value = 46382090
$ws.Range("I13").Value = 25365677
$ws.Range("J13").Value = 38224648
$ws.Range("K13").Value = 54548003
$ws.Range("L13").Value = 69105019
$ws.Range("M13").Value = 28110520
# Row 14
$ws.Range("D14").Value = -34001119
$ws.Range("E14").Value = -8851881
$ws.Range("F14").Value = -19517078
$ws.Range("G14").Value = -24276567
$ws.Range("H14").Value = -30463107
$ws.Range("I14").Value = -12827178
$ws.Range("J14").Value = -26593893
$ws.Range("K14").Value = -34963602
$ws.Range("L14").Value = -44188435
$ws.Range("M14").Value = -13527743
# Row 16
$ws.Range("D16").Value = 9770789
$ws.Range("E16").Value = -2118581
$ws.Range("F16").Value = -56678
$ws.Range("G16").Value = 152041
$ws.Range("H16").Value = -689512
$ws.Range("I16").Value = 971324
$ws.Range("J16").Value = 255139
$ws.Range("K16").Value = 4999540
$ws.Range("L16").Value = 5922920
$ws.Range("M16").Value = -2236562
# Row 17
$ws.Range("D17").Value = 57660214
$ws.Range("E17").Value = 19852229
$ws.Range("F17").Value = 27507833
$ws.Range("G17").Value = 16879120
$ws.Range("H17").Value = 15229471
$ws.Range("I17").Value = 13509823
$ws.Range("J17").Value = 11885894
$ws.Range("K17").Value = 24583941
$ws.Range("L17").Value = 30839504
$ws.Range("M17").Value = 12346215
# Row 19
$ws.Range("D19").Value = 247268
$ws.Range("E19").Value = 738279
$ws.Range("F19").Value = 1830909
$ws.Range("G19").Value = 2597851
$ws.Range("H19").Value = 3737347
$ws.Range("I19").Value = -2537743
$ws.Range("J19").Value = 1752854
$ws.Range("K19").Value = 2335410
$ws.Range("L19").Value = -479503
$ws.Range("M19").Value = -2302547
# Row 20
$ws.Range("D20").Value = 57907482
$ws.Range("E20").Value = 20590508
$ws.Range("F20").Value = 29338742
$ws.Range("G20").Value = 19476971
$ws.Range("H20").Value = 18966818
$ws.Range("I20").Value = 10972080
$ws.Range("J20").Value = 13638748
$ws.Range("K20").Value = 26919351
$ws.Range("L20").Value = 30360001
$ws.Range("M20").Value = 10043668
# Row 21
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = -5383870
$ws.Range("L21").Value = 0
# Row 22
$ws.Range("D22").Value = 57907482
$ws.Range("E22").Value = 20590508
$ws.Range("F22").Value = 29338742
$ws.Range("G22").Value = 19476971
$ws.Range("H22").Value = 18966818
$ws.Range("I22").Value = 10972080
$ws.Range("J22").Value = 13638748
$ws.Range("K22").Value = 21535481
$ws.Range("L22").Value = 30360001
$ws.Range("M22").Value = 10043668
# Row 24
$ws.Range("D24").Value = 57907482
$ws.Range("E24").Value = 20590508
$ws.Range("F24").Value = 29338742
$ws.Range("G24").Value = 19476971
$ws.Range("H24").Value = 18966818
$ws.Range("I24").Value = 10972080
$ws.Range("J24").Value = 13638748
$ws.Range("K24").Value = 21535481
$ws.Range("L24").Value = 30360001
$ws.Range("M24").Value = 10043668
# Row 25
$ws.Range("D25").Value = 24128
$ws.Range("E25").Value = 8579
$ws.Range("F25").Value = 12224
$ws.Range("G25").Value = 8115
$ws.Range("H25").Value = 7903
$ws.Range("I25").Value = 4572
$ws.Range("J25").Value = 5683
$ws.Range("K25").Value = 8973
$ws.Range("L25").Value = 12650
$ws.Range("M25").Value = 4185
# Row 27
$ws.Range("D27").Value = 24128
$ws.Range("E27").Value = 8579
$ws.Range("F27").Value = 12224
$ws.Range("G27").Value = 8115
$ws.Range("H27").Value = 7903
$ws.Range("I27").Value = 4572
$ws.Range("J27").Value = 5683
$ws.Range("K27").Value = 8973
$ws.Range("L27").Value = 12650
$ws.Range("M27").Value = 4185

# --- Column widths (shift by one quarter-column too) ---
$ws.Columns.Item(4).ColumnWidth = 28.17
$ws.Columns.Item(5).ColumnWidth = 27.17
$ws.Columns.Item(6).ColumnWidth = 27.17
$ws.Columns.Item(7).ColumnWidth = 27.17
$ws.Columns.Item(8).ColumnWidth = 28.17
$ws.Columns.Item(9).ColumnWidth = 27.17
$ws.Columns.Item(10).ColumnWidth = 27.17
$ws.Columns.Item(11).ColumnWidth = 27.17
$ws.Columns.Item(12).ColumnWidth = 28.17
$ws.Columns.Item(13).ColumnWidth = 27.17
